$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.25517891189561936
$ws.Range("B1").Value = 0.25490260027953582
$ws.Range("A2").Value = -0.18908241624958499
$ws.Range("B2").Value = 0.18839250726985757
$ws.Range("A3").Value = -0.13868381511172956
$ws.Range("B3").Value = 0.13845546900248884
$ws.Range("A4").Value = -0.13045546906325711
$ws.Range("B4").Value = 0.1297747464865715
$ws.Range("A5").Value = -0.12677474652382514
$ws.Range("B5").Value = 0.12443134850191928
$ws.Range("A6").Value = -0.025176391358700201
$ws.Range("B6").Value = 0.024856132675672882
$ws.Range("A7").Value = -0.014856132766137176
$ws.Range("B7").Value = 0.01478762145719692
$ws.Range("A8").Value = 0.0065143271619687582
$ws.Range("B8").Value = -0.0065293257928904502
$ws.Range("A9").Value = 0.008529325746759131
$ws.Range("B9").Value = -0.008537946494823867
$ws.Range("A10").Value = 0.010537946449266755
$ws.Range("B10").Value = -0.010537680698547547
$ws.Range("A11").Value = 0.013537680647355721
$ws.Range("B11").Value = -0.013538474244898957
$ws.Range("A12").Value = 0.017038474191288788
$ws.Range("B12").Value = -0.017064074261909035
$ws.Range("A13").Value = 0.020564074210523309
$ws.Range("B13").Value = -0.020591094154868372
$ws.Range("A14").Value = 0.028591094079094326
$ws.Range("B14").Value = -0.028640075395268028
$ws.Range("A15").Value = -0.0080516806115298323
$ws.Range("B15").Value = 0.0080336587146208061
$ws.Range("A16").Value = -0.0060336587556726329
$ws.Range("B16").Value = 0.0060033914393278209
$ws.Range("A17").Value = -0.0040033914813033533
$ws.Range("B17").Value = 0.0039999999465223368
$ws.Range("A18").Value = -0.059820184065689119
$ws.Range("B18").Value = 0.059733940752192183
$ws.Range("A19").Value = -0.012091661880651028
$ws.Range("B19").Value = 0.012016688936459907
$ws.Range("A20").Value = -0.0080166889635702177
$ws.Range("B20").Value = 0.0080056937858810073
$ws.Range("A21").Value = -0.004005693813286193
$ws.Range("B21").Value = 0.0039999999723647761
$ws.Range("A22").Value = -0.045708692191707812
$ws.Range("B22").Value = 0.04549642331686421
$ws.Range("A23").Value = -0.040496423357940792
$ws.Range("B23").Value = 0.040098496072144663
$ws.Range("A24").Value = -0.020098496201358174
$ws.Range("B24").Value = 0.019999999869092733
$ws.Range("A25").Value = -0.09725495718790178
$ws.Range("B25").Value = 0.097130207081253417
$ws.Range("A26").Value = -0.09463020713070236
$ws.Range("B26").Value = 0.094469388304650792
$ws.Range("A27").Value = -0.091969388356635484
$ws.Range("B27").Value = 0.091016045874316109
$ws.Range("A28").Value = -0.089016045934916299
$ws.Range("B28").Value = 0.088363357229118478
$ws.Range("A29").Value = -0.081363357325968444
$ws.Range("B29").Value = 0.08117417187944298
$ws.Range("A30").Value = -0.021174172278922043
$ws.Range("B30").Value = 0.021023682007331512
$ws.Range("A31").Value = -0.014023682111149682
$ws.Range("B31").Value = 0.014001255957994729
$ws.Range("A32").Value = -0.004001256079172677
$ws.Range("B32").Value = 0.0039999999124464836
